$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3" = 6529281
    "F3" = "FC Roskilde"
    "G3" = "Frem"
    "H3" = 3
    "I3" = 0
    "J3" = "H"
    "K3" = 1.85
    "L3" = 3.5
    "M3" = 3.75
    "N3" = 1.85
    "O3" = 3.4
    "P3" = 3.8
    "Q3" = -0.5
    "R3" = 1.825
    "S3" = 1.975
    "T3" = 2.5
    "U3" = 2
    "V3" = 1.8
    "W3" = 0.8500000000000001
    "X3" = -1
    "Y3" = -1
    "Z3" = 0.825
    "AA3" = -1
    "AB3" = 1
    "AC3" = -1
    "B5" = 6532917
    "F5" = "Aarhus Fremad"
    "G5" = "Kolding IF"
    "H5" = 1
    "I5" = 3
    "J5" = "A"
    "K5" = 2.7
    "L5" = 3.4
    "M5" = 2.35
    "N5" = 3.4
    "O5" = 3.6
    "P5" = 1.909
    "Q5" = 0.5
    "R5" = 1.85
    "S5" = 1.95
    "T5" = 2.25
    "U5" = 1.8
    "V5" = 2
    "W5" = -1
    "X5" = -1
    "Y5" = 0.909
    "Z5" = -1
    "AA5" = 0.95
    "AB5" = 0.8
    "AC5" = -1
    "B7" = 6529283
    "F7" = "Hellerup IK"
    "G7" = "Frem"
    "H7" = 4
    "I7" = 0
    "J7" = "H"
    "K7" = 1.95
    "L7" = 3.75
    "M7" = 3
    "N7" = 1.909
    "O7" = 3.8
    "P7" = 3.2
    "Q7" = -0.5
    "R7" = 1.95
    "S7" = 1.85
    "T7" = 3.25
    "U7" = 1.925
    "V7" = 1.875
    "W7" = 0.909
    "X7" = -1
    "Y7" = -1
    "Z7" = 0.95
    "AA7" = -1
    "AB7" = 0.925
    "AC7" = -1
    "B10" = 6532918
    "F10" = "AB Copenhagen"
    "G10" = "Kolding IF"
    "H10" = 1
    "I10" = 3
    "J10" = "A"
    "K10" = 5
    "L10" = 3.75
    "M10" = 1.571
    "N10" = 5
    "O10" = 3.75
    "P10" = 1.6
    "Q10" = 0.75
    "R10" = 2.05
    "S10" = 1.75
    "T10" = 2.5
    "U10" = 1.85
    "V10" = 1.95
    "W10" = -1
    "X10" = -1
    "Y10" = 0.6000000000000001
    "Z10" = -1
    "AA10" = 0.75
    "AB10" = 0.8500000000000001
    "AC10" = -1
    "B12" = 6532919
    "F12" = "Aarhus Fremad"
    "G12" = "AB Copenhagen"
    "H12" = 5
    "I12" = 2
    "J12" = "H"
    "K12" = 1.727
    "L12" = 3.8
    "M12" = 3.75
    "N12" = 1.5
    "O12" = 4.333
    "P12" = 5
    "Q12" = -1
    "R12" = 1.8
    "S12" = 2
    "T12" = 3.25
    "U12" = 1.95
    "V12" = 1.85
    "W12" = 0.5
    "X12" = -1
    "Y12" = -1
    "Z12" = 0.8
    "AA12" = -1
    "AB12" = 0.95
    "AC12" = -1
    "B14" = 6529284
    "F14" = "Frem"
    "G14" = "Brabrand"
    "H14" = 1
    "I14" = 1
    "J14" = "D"
    "K14" = 2.375
    "L14" = 3.8
    "M14" = 2.375
    "N14" = 2.4
    "O14" = 3.8
    "P14" = 2.3
    "Q14" = 0
    "R14" = 1.975
    "S14" = 1.825
    "T14" = 3
    "U14" = 1.95
    "V14" = 1.85
    "W14" = -1
    "X14" = 2.8
    "Y14" = -1
    "Z14" = 0
    "AA14" = 0
    "AB14" = -1
    "AC14" = 0.8500000000000001
    "B15" = 6533424
    "F15" = "Kolding IF"
    "G15" = "Esbjerg"
    "H15" = 1
    "I15" = 2
    "J15" = "A"
    "K15" = 2.5
    "L15" = 3.75
    "M15" = 2.3
    "N15" = 2.25
    "O15" = 4
    "P15" = 2.5
    "Q15" = 0
    "R15" = 1.775
    "S15" = 2.025
    "T15" = 2.75
    "U15" = 1.9
    "V15" = 1.9
    "W15" = -1
    "X15" = -1
    "Y15" = 1.5
    "Z15" = -1
    "AA15" = 1.025
    "AB15" = 0.45
    "AC15" = -0.5
    "B19" = 6858506
    "F19" = "FC Roskilde"
    "G19" = "FA 2000"
    "H19" = 1
    "I19" = 0
    "J19" = "H"
    "K19" = 1.571
    "L19" = 3.75
    "M19" = 4.75
    "N19" = 1.5
    "O19" = 3.8
    "P19" = 5.25
    "Q19" = -1
    "R19" = 1.9
    "S19" = 1.9
    "T19" = 2.75
    "U19" = 1.8
    "V19" = 2
    "W19" = 0.5
    "X19" = -1
    "Y19" = -1
    "Z19" = 0
    "AA19" = 0
    "AB19" = -1
    "AC19" = 1
    "B20" = 6858902
    "F20" = "Brabrand"
    "G20" = "Thisted FC"
    "H20" = 1
    "I20" = 2
    "J20" = "A"
    "K20" = 2.3
    "L20" = 3.3
    "M20" = 2.625
    "N20" = 2.15
    "O20" = 3.3
    "P20" = 2.8
    "Q20" = -0.25
    "R20" = 1.95
    "S20" = 1.85
    "T20" = 2.5
    "U20" = 1.875
    "V20" = 1.925
    "W20" = -1
    "X20" = -1
    "Y20" = 1.8
    "Z20" = -1
    "AA20" = 0.8500000000000001
    "AB20" = 0.875
    "AC20" = -1
    "B29" = 6858911
    "F29" = "Skive"
    "G29" = "Aarhus Fremad"
    "H29" = 1
    "I29" = 2
    "J29" = "A"
    "K29" = 4.75
    "L29" = 3.8
    "M29" = 1.6
    "N29" = 5
    "O29" = 3.8
    "P29" = 1.571
    "Q29" = 1
    "R29" = 1.8
    "S29" = 2
    "T29" = 2.75
    "U29" = 1.75
    "V29" = 1.95
    "W29" = -1
    "X29" = -1
    "Y29" = 0.571
    "Z29" = 0
    "AA29" = 0
    "AB29" = 0.375
    "AC29" = -0.5
    "B30" = 6858913
    "F30" = "Nykobing"
    "G30" = "FC Roskilde"
    "H30" = 1
    "I30" = 2
    "J30" = "A"
    "K30" = 2.05
    "L30" = 3.5
    "M30" = 3.1
    "N30" = 2.15
    "O30" = 3.4
    "P30" = 2.9
    "Q30" = -0.25
    "R30" = 1.925
    "S30" = 1.875
    "T30" = 2.75
    "U30" = 1.875
    "V30" = 1.925
    "W30" = -1
    "X30" = -1
    "Y30" = 1.9
    "Z30" = -1
    "AA30" = 0.875
    "AB30" = 0.4375
    "AC30" = -0.5
    "B82" = 6859007
    "F82" = "Skive"
    "G82" = "FC Roskilde"
    "H82" = 1
    "I82" = 2
    "J82" = "A"
    "K82" = 3.6
    "L82" = 3.4
    "M82" = 1.909
    "N82" = 3.2
    "O82" = 3.4
    "P82" = 2.05
    "Q82" = 0.25
    "R82" = 2
    "S82" = 1.8
    "T82" = 2.75
    "U82" = 1.975
    "V82" = 1.825
    "W82" = -1
    "X82" = -1
    "Y82" = 1.05
    "Z82" = -1
    "AA82" = 0.8
    "AB82" = 0.4875
    "AC82" = -0.5
    "B83" = 6859008
    "F83" = "Brabrand"
    "G83" = "AB Copenhagen"
    "H83" = 2
    "I83" = 2
    "J83" = "D"
    "K83" = 3.6
    "L83" = 3.6
    "M83" = 1.85
    "N83" = 3.4
    "O83" = 3.6
    "P83" = 1.909
    "Q83" = 0.5
    "R83" = 1.825
    "S83" = 1.975
    "T83" = 2.5
    "U83" = 1.8
    "V83" = 2
    "W83" = -1
    "X83" = 2.6
    "Y83" = -1
    "Z83" = 0.825
    "AA83" = -1
    "AB83" = 0.8
    "AC83" = -1
    "B85" = 6859010
    "F85" = "Esbjerg"
    "G85" = "FA 2000"
    "H85" = 3
    "I85" = 1
    "J85" = "H"
    "K85" = 1.222
    "L85" = 6.5
    "M85" = 9
    "N85" = 1.2
    "O85" = 6.5
    "P85" = 10
    "Q85" = -2
    "R85" = 1.9
    "S85" = 1.9
    "T85" = 3.5
    "U85" = 1.875
    "V85" = 1.925
    "W85" = 0.2
    "X85" = -1
    "Y85" = -1
    "Z85" = 0
    "AA85" = 0
    "AB85" = 0.875
    "AC85" = -1
    "B118" = 6859068
    "F118" = "Fremad Amager"
    "G118" = "FA 2000"
    "H118" = 0
    "I118" = 0
    "J118" = "D"
    "K118" = 2.2
    "L118" = 3.4
    "M118" = 2.9
    "N118" = 2.15
    "O118" = 3.4
    "P118" = 2.9
    "Q118" = -0.25
    "R118" = 1.925
    "S118" = 1.875
    "T118" = 2.5
    "U118" = 1.925
    "V118" = 1.875
    "W118" = -1
    "X118" = 2.4
    "Y118" = -1
    "Z118" = -0.5
    "AA118" = 0.4375
    "AB118" = -1
    "AC118" = 0.875
    "B119" = 6859066
    "F119" = "Hellerup IK"
    "G119" = "FC Roskilde"
    "H119" = 2
    "I119" = 3
    "J119" = "A"
    "K119" = 3.6
    "L119" = 3.6
    "M119" = 1.833
    "N119" = 4.5
    "O119" = 3.8
    "P119" = 1.615
    "Q119" = 0.75
    "R119" = 1.975
    "S119" = 1.825
    "T119" = 2.75
    "U119" = 1.875
    "V119" = 1.925
    "W119" = -1
    "X119" = -1
    "Y119" = 0.615
    "Z119" = -0.5
    "AA119" = 0.4125
    "AB119" = 0.875
    "AC119" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}